# Update the "想去人数" (want-to-go count) column F values
# across all four worksheets to match the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 218
$ws.Range("F4").Value  = 401
$ws.Range("F5").Value  = 201
$ws.Range("F6").Value  = 803
$ws.Range("F7").Value  = 103
$ws.Range("F8").Value  = 10264
$ws.Range("F10").Value = 3537
$ws.Range("F12").Value = 2452
$ws.Range("F13").Value = 35
$ws.Range("F14").Value = 2825
$ws.Range("F17").Value = 2184
$ws.Range("F19").Value = 97
$ws.Range("F20").Value = 27
$ws.Range("F23").Value = 153
$ws.Range("F24").Value = 318
$ws.Range("F25").Value = 276
$ws.Range("F26").Value = 231
$ws.Range("F27").Value = 617
$ws.Range("F28").Value = 1322
$ws.Range("F30").Value = 1259
$ws.Range("F31").Value = 106
$ws.Range("F34").Value = 3786
$ws.Range("F35").Value = 3200
$ws.Range("F36").Value = 33
$ws.Range("F38").Value = 1045
$ws.Range("F39").Value = 402
$ws.Range("F41").Value = 1294
$ws.Range("F42").Value = 104
$ws.Range("F43").Value = 111
$ws.Range("F47").Value = 12

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 9
$ws.Range("F16").Value = 180

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 754
$ws.Range("F4").Value = 128
$ws.Range("F5").Value = 2044

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 754
$ws.Range("F4").Value  = 128
$ws.Range("F6").Value  = 401
$ws.Range("F8").Value  = 201
$ws.Range("F9").Value  = 803
$ws.Range("F10").Value = 103
$ws.Range("F11").Value = 10264
$ws.Range("F13").Value = 3537
$ws.Range("F15").Value = 2452
$ws.Range("F16").Value = 35
$ws.Range("F19").Value = 2184
$ws.Range("F21").Value = 97
$ws.Range("F22").Value = 27
$ws.Range("F24").Value = 153
$ws.Range("F25").Value = 318
$ws.Range("F26").Value = 231
$ws.Range("F27").Value = 1322
$ws.Range("F29").Value = 1259
$ws.Range("F30").Value = 106
$ws.Range("F33").Value = 9
$ws.Range("F36").Value = 3200
$ws.Range("F37").Value = 1045
$ws.Range("F44").Value = 1294
$ws.Range("F45").Value = 104
$ws.Range("F48").Value = 12
$ws.Range("F49").Value = 180

$wb.Save()
